{"js": "// Commit: \"rtp test in unity\"\n//\n// Appends a new \"13. h\u00e9t\" weekly-report entry at the end of the document,\n// matching the justified alignment + Hungarian language formatting used by\n// every other paragraph in the log:\n//   1. \"13. h\u00e9t\"                                   (heading line)\n//   2. \"A WebRTC hossz\u00fa sikertelens\u00e9ge ...\"         (first body paragraph)\n//   3. \"Az RTP m\u0171k\u00f6d\u00e9s\u00e9nek m\u00e1snapj\u00e1n ...\"           (second body paragraph)\n//   4. \"A VLC Unity al\u00e1 k\u00e9sz\u00fclt libVLC ...\"         (third body paragraph)\n\nconst body = context.document.body;\n\nconst newParagraphs = [\n  {\n    text: \"13. h\u00e9t\",\n    lastRenderedPageBreak: false\n  },\n  {\n    text: \"A WebRTC hossz\u00fa sikertelens\u00e9ge ut\u00e1n \u00e1tt\u00e9r\u00e9st kezdtem RTP-re, amit viszonylag gyorsan siker\u00fclt m\u0171k\u00f6d\u00e9sre b\u00edrnom a virtu\u00e1lis ubuntu \u00e9s a gazdag\u00e9pen fut\u00f3 VLC k\u00f6z\u00f6tt. Mivel a VLC-nek van Unity pluginje, ez\u00e9rt bizakod\u00f3 vagyok.\",\n    lastRenderedPageBreak: true\n  },\n  {\n    text: \"Az RTP m\u0171k\u00f6d\u00e9s\u00e9nek m\u00e1snapj\u00e1n reggel m\u00e1r nem m\u0171k\u00f6d\u00f6tt a dolog, illetve egy virtu\u00e1lis g\u00e9p restart ut\u00e1n nem jelent meg semmilyen interakt\u00e1lhat\u00f3 grafikai elem az Ubuntu UI-on, amit eg\u00e9sz probl\u00e9m\u00e1s volt visszaszerezni.\",\n    lastRenderedPageBreak: false\n  },\n  {\n    text: \"A VLC Unity al\u00e1 k\u00e9sz\u00fclt libVLC implement\u00e1ci\u00f3j\u00e1nak seg\u00edts\u00e9g\u00e9vel siker\u00fclt rtp streamet megjelen\u00edteni a virtu\u00e1lis g\u00e9pr\u0151l.\",\n    lastRenderedPageBreak: false\n  }\n];\n\nfor (const info of newParagraphs) {\n  // Insert the new paragraph at the end of the body and give it the same\n  // \"justified + hu-HU\" formatting every other paragraph in this log uses.\n  const paragraph = body.insertParagraph(info.text, \"End\");\n  paragraph.alignment = Word.Alignment.justified;\n  paragraph.font.languageId = Word.LanguageId.hungarian;\n\n  if (info.lastRenderedPageBreak) {\n    // Word stamps a <w:lastRenderedPageBreak/> marker at the point where a\n    // page boundary fell the last time the document was laid out/rendered.\n    // Reproduce that exact marker (immediately before the run text) via a\n    // small raw-OOXML replace of the paragraph we just created.\n    await context.sync();\n    const range = paragraph.getRange(\"Whole\");\n    const ooxml =\n      '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n      '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n      '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n      '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n      '<w:body><w:p><w:pPr><w:jc w:val=\"both\"/><w:rPr><w:lang w:val=\"hu-HU\"/></w:rPr></w:pPr>' +\n      '<w:r><w:rPr><w:lang w:val=\"hu-HU\"/></w:rPr><w:lastRenderedPageBreak/>' +\n      '<w:t>' + info.text + '</w:t></w:r></w:p></w:body>' +\n      '</w:document></pkg:xmlData></pkg:part></pkg:package>';\n    range.insertOoxml(ooxml, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Commit: \"rtp test in unity\"\n#\n# Appends a new \"13. h\u00e9t\" weekly-report entry at the end of the document,\n# matching the justified alignment + Hungarian language formatting every\n# other paragraph in this log already uses:\n#   1. \"13. h\u00e9t\"                                   (heading line)\n#   2. \"A WebRTC hossz\u00fa sikertelens\u00e9ge ...\"         (first body paragraph)\n#   3. \"Az RTP m\u0171k\u00f6d\u00e9s\u00e9nek m\u00e1snapj\u00e1n ...\"           (second body paragraph)\n#   4. \"A VLC Unity al\u00e1 k\u00e9sz\u00fclt libVLC ...\"         (third body paragraph)\n\n$d = $word.ActiveDocument\n\n$wdAlignParagraphJustify = 3\n\n$newParagraphs = @(\n  @{ Text = \"13. h\u00e9t\"; PageBreak = $false },\n  @{ Text = \"A WebRTC hossz\u00fa sikertelens\u00e9ge ut\u00e1n \u00e1tt\u00e9r\u00e9st kezdtem RTP-re, amit viszonylag gyorsan siker\u00fclt m\u0171k\u00f6d\u00e9sre b\u00edrnom a virtu\u00e1lis ubuntu \u00e9s a gazdag\u00e9pen fut\u00f3 VLC k\u00f6z\u00f6tt. Mivel a VLC-nek van Unity pluginje, ez\u00e9rt bizakod\u00f3 vagyok.\"; PageBreak = $true },\n  @{ Text = \"Az RTP m\u0171k\u00f6d\u00e9s\u00e9nek m\u00e1snapj\u00e1n reggel m\u00e1r nem m\u0171k\u00f6d\u00f6tt a dolog, illetve egy virtu\u00e1lis g\u00e9p restart ut\u00e1n nem jelent meg semmilyen interakt\u00e1lhat\u00f3 grafikai elem az Ubuntu UI-on, amit eg\u00e9sz probl\u00e9m\u00e1s volt visszaszerezni.\"; PageBreak = $false },\n  @{ Text = \"A VLC Unity al\u00e1 k\u00e9sz\u00fclt libVLC implement\u00e1ci\u00f3j\u00e1nak seg\u00edts\u00e9g\u00e9vel siker\u00fclt rtp streamet megjelen\u00edteni a virtu\u00e1lis g\u00e9pr\u0151l.\"; PageBreak = $false }\n)\n\nforeach ($item in $newParagraphs) {\n    # Insert a fresh paragraph after the current last one; it inherits the\n    # \"justified + hu-HU\" formatting already used throughout the document,\n    # but we also set both explicitly to be safe.\n    $lastRange = $d.Paragraphs.Last.Range\n    $lastRange.InsertParagraphAfter()\n\n    $newPara = $d.Paragraphs.Last\n    $newPara.Range.Text = $item.Text\n    $newPara.Alignment = $wdAlignParagraphJustify\n    $newPara.Range.LanguageID = \"hu-HU\"\n\n    if ($item.PageBreak) {\n        # Word stamps a <w:lastRenderedPageBreak/> marker at the spot where a\n        # page boundary fell the last time the document was laid out. Re-create\n        # that exact marker (immediately before the run text) via a scoped\n        # InsertXML replace of the paragraph we just created.\n        $wholeRange = $newPara.Range\n        $xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:jc w:val=\"both\"/><w:rPr><w:lang w:val=\"hu-HU\"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=\"hu-HU\"/></w:rPr><w:lastRenderedPageBreak/><w:t>' + $item.Text + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n        $wholeRange.InsertXML($xml)\n    }\n}\n"}
